# Append one new data row (row 21) to the NIFTY_Options_Analysis sheet,
# mirroring the layout/formatting of the last existing row (row 20) with
# this cycle's values (AVOID / AVOID, bearish CPR trending-day veto).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 21
$styleRow = 20

# ---------------------------------------------------------------------
# 1) Stamp the new row's cells with the same visual formatting as the
#    template row (center/middle aligned, word-wrap, thin box border)
#    before writing any values - this lets the engine fold each cell
#    back onto the existing shared style instead of minting new ones.
# ---------------------------------------------------------------------
function Format-LikeRow20($col) {
    $cell = $ws.Cells.Item($newRow, $col)
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4108     # xlCenter
    $cell.WrapText = $true
    $cell.Borders.LineStyle = 1         # thin box border
    return $cell
}

# Plain text columns (style matches row 20's default text style)
$textCols = 1,2,5,6,13,15,20,21,29,30,31   # A,B,E,F,M,O,T,U,AC,AD,AE
foreach ($col in $textCols) {
    [void](Format-LikeRow20 $col)
}

# Numeric columns with "0.0" format (G,K,L,N,P,Q,R,S)
foreach ($col in 7,11,12,14,16,17,18,19) {
    $cell = Format-LikeRow20 $col
    $cell.NumberFormat = "0.0"
}

# Numeric columns with "0.00" format (H,I,J,W,X,Z,AA)
foreach ($col in 8,9,10,23,24,26,27) {
    $cell = Format-LikeRow20 $col
    $cell.NumberFormat = "0.00"
}

# Numeric columns with "0.000000" format (Y,AB)
foreach ($col in 25,28) {
    $cell = Format-LikeRow20 $col
    $cell.NumberFormat = "0.000000"
}

# V column: plain/default-styled but holds a numeric value
[void](Format-LikeRow20 22)

# C (Signal=AVOID) and D (Signal_Tier=AVOID) carry the bold
# warning/alert look used throughout the sheet for AVOID rows.
$cellC = Format-LikeRow20 3
$cellC.Font.Bold = $true
$cellC.Font.Color = 0x06009C        # BGR for RGB 9C0006 (dark red text)
$cellC.Interior.Color = 0xCEC7FF    # BGR for RGB FFC7CE (light red fill)

$cellD = Format-LikeRow20 4
$cellD.Font.Bold = $true
$cellD.Font.Color = 0xFFFFFF        # white text
$cellD.Interior.Color = 0x0000FF    # BGR for RGB FF0000 (solid red fill)

# ---------------------------------------------------------------------
# 2) Write the row's values.
#    A21 ("2026-01-19") and E21 ("100%") look like a date / a percentage
#    to Excel's input parser, so they are entered with a leading
#    apostrophe to force literal text, matching the stored text in the
#    source row. Every other text cell is plain text already.
# ---------------------------------------------------------------------
$ws.Cells.Item($newRow, 1).Value  = "'2026-01-19"            # A  Date
$ws.Cells.Item($newRow, 2).Value  = "10:00:08"                # B  Time
$ws.Cells.Item($newRow, 3).Value  = "AVOID"                   # C  Signal
$ws.Cells.Item($newRow, 4).Value  = "AVOID"                   # D  Signal_Tier
$ws.Cells.Item($newRow, 5).Value  = "'100%"                   # E  Position_Size
$ws.Cells.Item($newRow, 6).Value  = "TRADEABLE"                # F  Premium_Quality
$ws.Cells.Item($newRow, 7).Value  = 0                          # G  Total_Score
$ws.Cells.Item($newRow, 8).Value  = 25539.9                    # H  NIFTY_Spot
$ws.Cells.Item($newRow, 9).Value  = 11.98                      # I  VIX
$ws.Cells.Item($newRow, 10).Value = 0.78                       # J  VIX_Trend
$ws.Cells.Item($newRow, 11).Value = 0                          # K  VIX_Score
$ws.Cells.Item($newRow, 12).Value = 39.6                       # L  IV_Rank
$ws.Cells.Item($newRow, 13).Value = "UNKNOWN"                   # M  Market_Regime
$ws.Cells.Item($newRow, 14).Value = 0                           # N  Regime_Score
$ws.Cells.Item($newRow, 15).Value = "UNKNOWN"                   # O  OI_Pattern
$ws.Cells.Item($newRow, 16).Value = 0                           # P  OI_Score
$ws.Cells.Item($newRow, 17).Value = 0                           # Q  Theta_Score
$ws.Cells.Item($newRow, 18).Value = 0                           # R  Gamma_Score
$ws.Cells.Item($newRow, 19).Value = 0                           # S  Vega_Score
$ws.Cells.Item($newRow, 20).Value = "NONE"                      # T  Best_Strategy
$ws.Cells.Item($newRow, 21).Value = ""                          # U  Expiry_1 (blank)
$ws.Cells.Item($newRow, 22).Value = 0                           # V  Days_To_Expiry_1
$ws.Cells.Item($newRow, 23).Value = 0                           # W  Straddle_Premium
$ws.Cells.Item($newRow, 24).Value = 0                           # X  Straddle_Theta
$ws.Cells.Item($newRow, 25).Value = 0                           # Y  Straddle_Gamma
$ws.Cells.Item($newRow, 26).Value = 0                           # Z  Strangle_Premium
$ws.Cells.Item($newRow, 27).Value = 0                           # AA Strangle_Theta
$ws.Cells.Item($newRow, 28).Value = 0                           # AB Strangle_Gamma
$ws.Cells.Item($newRow, 29).Value = "HARD VETO: CPR TRENDING DAY: Price 25539.90 below BC 25767.95 - BEARISH TRENDING DAY likely"   # AC Recommendation
$ws.Cells.Item($newRow, 30).Value = "CPR TRENDING DAY: Price 25539.90 below BC 25767.95 - BEARISH TRENDING DAY likely"              # AD Risk_Factors
$ws.Cells.Item($newRow, 31).Value = "Yes"                       # AE Telegram_Sent

Write-Output "Row 21 written"
